$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the custom column width on column A (cols element disappears in the diff)
$ws.Columns.Item(1).AutoFit() | Out-Null

# Row 2: Total OTUs
$ws.Range("B2").Value = 268
$ws.Range("C2").Value = 202
$ws.Range("D2").Value = 421
$ws.Range("E2").Value = 364

# Row 3: Unique OTUs
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 64
$ws.Range("E3").Value = 212

# Row 4: Shared with Epiphytes
$ws.Range("C4").Value = 180
$ws.Range("D4").Value = 247
$ws.Range("E4").Value = 57

# Row 5: Shared with Endophytes
$ws.Range("B5").Value = 180
$ws.Range("D5").Value = 187
$ws.Range("E5").Value = 43

# Row 6: Shared with Litter
$ws.Range("B6").Value = 247
$ws.Range("C6").Value = 187
$ws.Range("E6").Value = 145

# Row 7: Total Shared (percentages)
$ws.Range("B7").Value = "263 (98%)"
$ws.Range("C7").Value = "200 (99%)"
$ws.Range("D7").Value = "357 (85%)"
$ws.Range("E7").Value = "152 (42%)"

# Selection matching the diff
$ws.Range("H7").Select() | Out-Null
